$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B12").ClearContents()
$ws.Range("B13:F13").Cut($ws.Range("B11:F11"))
$ws.Range("B13:F13").Clear()
$ws.Range("B11:F11").Select()
